$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the two brand-new leaderboard rows first (shifts old rows down) ---
$ws.Rows.Item(7).Insert()
# after this insert, rows are: 7=blank(new) 8=old7 ... 13=old12(236858) 14=old13(484487)
$ws.Rows.Item(15).Insert()
# after this insert, rows are: 7=blank(new) 8=old7 ... 14=old13(236858) 15=blank(new) 16=old14(484487)

# --- Populate the new row 15 first so its new shared string (21f1002538) is registered before the others ---
$ws.Range("B15").Value = "https://leetcode.com/u/21f1002538/"
$ws.Range("A15").Value = 254342
$ws.Range("A15").NumberFormat = "#,##0"
$ws.Range("D15").Value = 254
$ws.Range("E15").Value = 31
$ws.Range("F15").Value = 19

# --- Populate the new row 7 ---
$ws.Range("B7").Value = "https://leetcode.com/u/MikPosp/"
# A7 needs the same cell style as the rank cells above it (e.g. A6) - copy format then overwrite value
$ws.Range("A6").Copy()
$ws.Range("A7").PasteSpecial(-4104)
$excel.CutCopyMode = 0
$ws.Range("A7").Value = 2602
$ws.Range("C7").Value = 15
$ws.Range("D7").Value = 1420
$ws.Range("F7").Value = 77
$ws.Range("L7").Value = "No data"

# --- Restore the selection to match the saved view state ---
$ws.Range("L7").Select()
